# Rename the first sheet ("Spouse 1" -> "Sam") and make it the active
# (selected) sheet, replacing "Spouse 2" as the active tab.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Name = "Sam"

# Activating sheet1 makes it the workbook's active tab and marks its
# sheetView as the selected one, clearing the flag on the previously
# active sheet (Spouse 2).
$sheet1.Activate()
